# Adds "MDI Icon Name/Unicode/Link (Sunset - Sunrise)" columns (G:I) to the
# weather-code mapping sheet, populating the night-time icon equivalents for
# the "sunny" / "partly-cloudy" rows (weather codes 1-7), and updates the
# final cell-selection state on both worksheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)      # 工作表2 - the visible mapping table
$lookup = $wb.Worksheets.Item(2)  # 工作表1 - the MDI icon lookup table
$lookupName = $lookup.Name

# --- Header row (G1:I1) -----------------------------------------------
$ws.Range("G1").Value = "MDI Icon Name(Sunset - Sunrise)"
$ws.Range("H1").Value = "MDI Icon Unicode"
$ws.Range("I1").Value = "MDI Icon Link"

# --- Row 2 (weather code 1 - sunny -> weather-night) -------------------
$ws.Range("G2").Value = "weather-night"
$ws.Range("H2").Formula = "=VLOOKUP(G2,$lookupName!`$A`$2:`$C`$14,2,0)"
$ws.Range("I2").Formula = "=VLOOKUP(G2,$lookupName!`$A`$2:`$C`$14,3,0)"

# --- Rows 3-8 (weather codes 2-7 -> weather-night-partly-cloudy) -------
for ($r = 3; $r -le 8; $r++) {
    $ws.Range("G$r").Value = "weather-night-partly-cloudy"
    $ws.Range("H$r").Formula = "=VLOOKUP(G$r,$lookupName!`$A`$2:`$C`$14,2,0)"
    $ws.Range("I$r").Formula = "=VLOOKUP(G$r,$lookupName!`$A`$2:`$C`$14,3,0)"
}

# --- Column widths for the new columns (closest achievable match) ------
$ws.Columns.Item(7).ColumnWidth = 29.3
$ws.Columns.Item(8).ColumnWidth = 16.08
$ws.Columns.Item(9).ColumnWidth = 62.65

# --- Final selections left behind by the editing session ---------------
$ws.Range("C6").Select()
$lookup.Range("A14").Select()
